$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 15.67081492403954
$ws.Range("C2").Value = 4.920529886061695
$ws.Range("D2").Value = 7.999223543718799
$ws.Range("E2").Value = 10.171690038876
$ws.Range("F2").Value = 39.13785774348465
$ws.Range("I2").Value = 32.23090148857825
$ws.Range("K2").Value = 13.68377707550449
$ws.Range("L2").Value = 10.40803754457155
$ws.Range("M2").Value = 15.90298431543088
$ws.Range("N2").Value = 23.27016272231283

# Row 3
$ws.Range("B3").Value = 15.51084573806486
$ws.Range("C3").Value = 4.652768619461256
$ws.Range("D3").Value = 8.004217635475538
$ws.Range("E3").Value = 10.18728975260468
$ws.Range("F3").Value = 39.08171923743665
$ws.Range("I3").Value = 32.26242961528229
$ws.Range("K3").Value = 13.56592547323185
$ws.Range("L3").Value = 10.41772352254618
$ws.Range("M3").Value = 15.89026421742997
$ws.Range("N3").Value = 23.32136920468548

# Row 4
$ws.Range("B4").Value = 15.41588974984381
$ws.Range("C4").Value = 4.479307964505526
$ws.Range("D4").Value = 8.007336526809727
$ws.Range("E4").Value = 10.19768113381811
$ws.Range("F4").Value = 39.05552495011164
$ws.Range("I4").Value = 32.28709074116808
$ws.Range("K4").Value = 13.49652669539644
$ws.Range("L4").Value = 10.42513447020948
$ws.Range("M4").Value = 15.88529043399297
$ws.Range("N4").Value = 23.35470868525066

# Row 5
$ws.Range("B5").Value = 15.37805888224996
$ws.Range("C5").Value = 4.406370221192251
$ws.Range("D5").Value = 8.008620847805497
$ws.Range("E5").Value = 10.20212056106499
$ws.Range("F5").Value = 39.04693701677872
$ws.Range("I5").Value = 32.29847254409292
$ws.Range("K5").Value = 13.46902004518511
$ws.Range("L5").Value = 10.42852285962408
$ws.Range("M5").Value = 15.88397935991169
$ws.Range("N5").Value = 23.36877271179199

# Row 6
$ws.Range("B6").Value = 15.3718305368675
$ws.Range("C6").Value = 4.394123872829045
$ws.Range("D6").Value = 8.008834918770487
$ws.Range("E6").Value = 10.20287010893349
$ws.Range("F6").Value = 39.04563714665295
$ws.Range("I6").Value = 32.30044290319895
$ws.Range("K6").Value = 13.46450014566092
$ws.Range("L6").Value = 10.42910775522105
$ws.Range("M6").Value = 15.88380495520749
$ws.Range("N6").Value = 23.37113691341548

# Row 7
$ws.Range("B7").Value = 15.41537599349645
$ws.Range("C7").Value = 4.47833337693949
$ws.Range("D7").Value = 8.007353793388146
$ws.Range("E7").Value = 10.19774017550417
$ws.Range("F7").Value = 39.05540067611001
$ws.Range("I7").Value = 32.28723884804805
$ws.Range("K7").Value = 13.49615256043135
$ws.Range("L7").Value = 10.42517867535329
$ws.Range("M7").Value = 15.88526985115366
$ws.Range("N7").Value = 23.35489642189019

# Row 8
$ws.Range("B8").Value = 15.61500622883878
$ws.Range("C8").Value = 4.830098539612255
$ws.Range("D8").Value = 8.0009346990607
$ws.Range("E8").Value = 10.17690030058457
$ws.Range("F8").Value = 39.1167875664426
$ws.Range("I8").Value = 32.24067115670103
$ws.Range("K8").Value = 13.64254477777401
$ws.Range("L8").Value = 10.41107367539833
$ws.Range("M8").Value = 15.89801136932094
$ws.Range("N8").Value = 23.28742500113151

# Row 9
$ws.Range("B9").Value = 16.03041317888718
$ws.Range("C9").Value = 5.447225721337789
$ws.Range("D9").Value = 7.988756511359727
$ws.Range("E9").Value = 10.14246741738069
$ws.Range("F9").Value = 39.30251150559771
$ws.Range("I9").Value = 32.19147216718071
$ws.Range("K9").Value = 13.95178599548904
$ws.Range("L9").Value = 10.39501275510289
$ws.Range("M9").Value = 15.94537612566926
$ws.Range("N9").Value = 23.17015333071658

# Row 10
$ws.Range("B10").Value = 16.34746526711157
$ws.Range("C10").Value = 5.855514113362812
$ws.Range("D10").Value = 7.980048883276543
$ws.Range("E10").Value = 10.12106829534417
$ws.Range("F10").Value = 39.4782820849949
$ws.Range("I10").Value = 32.1810514585062
$ws.Range("K10").Value = 14.19062155000902
$ws.Range("L10").Value = 10.39026118279254
$ws.Range("M10").Value = 15.99362533723012
$ws.Range("N10").Value = 23.09312848792201

# Row 11
$ws.Range("B11").Value = 16.49367675634403
$ws.Range("C11").Value = 6.03136069479659
$ws.Range("D11").Value = 7.976137439966094
$ws.Range("E11").Value = 10.11217483216917
$ws.Range("F11").Value = 39.56663885232348
$ws.Range("I11").Value = 32.18190047364126
$ws.Range("K11").Value = 14.30138320636042
$ws.Range("L11").Value = 10.38962382720508
$ws.Range("M11").Value = 16.01844655295976
$ws.Range("N11").Value = 23.06006447103059

# Row 12
$ws.Range("B12").Value = 16.54927801328272
$ws.Range("C12").Value = 6.096522621389722
$ws.Range("D12").Value = 7.974663265944256
$ws.Range("E12").Value = 10.10892765713568
$ws.Range("F12").Value = 39.60128997984646
$ws.Range("I12").Value = 32.18302528627876
$ws.Range("K12").Value = 14.34359398046486
$ws.Range("L12").Value = 10.3896009694151
$ws.Range("M12").Value = 16.02825388437267
$ws.Range("N12").Value = 23.04782753356896

# Row 13
$ws.Range("B13").Value = 16.53729368363027
$ws.Range("C13").Value = 6.082552386139155
$ws.Range("D13").Value = 7.974980446411991
$ws.Range("E13").Value = 10.10962163705696
$ws.Range("F13").Value = 39.59377446976414
$ws.Range("I13").Value = 32.18274732152745
$ws.Range("K13").Value = 14.33449180699859
$ws.Range("L13").Value = 10.38959618626667
$ws.Range("M13").Value = 16.02612363420489
$ws.Range("N13").Value = 23.050450364621

# Row 14
$ws.Range("B14").Value = 16.49824670394986
$ws.Range("C14").Value = 6.036750237594539
$ws.Range("D14").Value = 7.976016019210023
$ws.Range("E14").Value = 10.11190527047502
$ws.Range("F14").Value = 39.56946581824906
$ws.Range("I14").Value = 32.18197691706338
$ws.Range("K14").Value = 14.3048507547928
$ws.Range("L14").Value = 10.38961757188133
$ws.Range("M14").Value = 16.01924525946959
$ws.Range("N14").Value = 23.05905204699991

# Row 15
$ws.Range("B15").Value = 16.47435828370176
$ws.Range("C15").Value = 6.008509080215772
$ws.Range("D15").Value = 7.976651245459429
$ws.Range("E15").Value = 10.11331975500767
$ws.Range("F15").Value = 39.55473086490601
$ws.Range("I15").Value = 32.18160961685633
$ws.Range("K15").Value = 14.28672854159024
$ws.Range("L15").Value = 10.38965910420987
$ws.Range("M15").Value = 16.01508504545792
$ws.Range("N15").Value = 23.06435776150712

# Row 16
$ws.Range("B16").Value = 16.33794553951184
$ws.Range("C16").Value = 5.843822652250121
$ws.Range("D16").Value = 7.980305493217377
$ws.Range("E16").Value = 10.12166639820328
$ws.Range("F16").Value = 39.47267548008034
$ws.Range("I16").Value = 32.18110841784942
$ws.Range("K16").Value = 14.18342246244547
$ws.Range("L16").Value = 10.39033344938797
$ws.Range("M16").Value = 15.99206059660217
$ws.Range("N16").Value = 23.09532901369097

# Row 17
$ws.Range("B17").Value = 16.25473227284353
$ws.Range("C17").Value = 5.740256028934024
$ws.Range("D17").Value = 7.982559883858029
$ws.Range("E17").Value = 10.12700196212962
$ws.Range("F17").Value = 39.42447769505469
$ws.Range("I17").Value = 32.18223234306669
$ws.Range("K17").Value = 14.12056329490608
$ws.Range("L17").Value = 10.39113705548434
$ws.Range("M17").Value = 15.9786681243238
$ws.Range("N17").Value = 23.11483448443402

# Row 18
$ws.Range("B18").Value = 16.20705989357783
$ws.Range("C18").Value = 5.679756932961864
$ws.Range("D18").Value = 7.983861236152974
$ws.Range("E18").Value = 10.13015002791391
$ws.Range("F18").Value = 39.39754706446407
$ws.Range("I18").Value = 32.18340495094236
$ws.Range("K18").Value = 14.08460955379594
$ws.Range("L18").Value = 10.39174277504131
$ws.Range("M18").Value = 15.97123581623384
$ws.Range("N18").Value = 23.12623942610151

# Row 19
$ws.Range("B19").Value = 16.19095298147171
$ws.Range("C19").Value = 5.659113409033314
$ws.Range("D19").Value = 7.98430266094054
$ws.Range("E19").Value = 10.13122952100764
$ws.Range("F19").Value = 39.38856518894162
$ws.Range("I19").Value = 32.18389235104104
$ws.Range("K19").Value = 14.07247189031445
$ws.Range("L19").Value = 10.39197252799232
$ws.Range("M19").Value = 15.96876600028642
$ws.Range("N19").Value = 23.13013288354628

# Row 20
$ws.Range("B20").Value = 16.26357120492086
$ws.Range("C20").Value = 5.751377205552342
$ws.Range("D20").Value = 7.982319416118335
$ws.Range("E20").Value = 10.12642578895521
$ws.Range("F20").Value = 39.42952662381807
$ws.Range("I20").Value = 32.18205824612653
$ws.Range("K20").Value = 14.12723420023628
$ws.Range("L20").Value = 10.39103666250125
$ws.Range("M20").Value = 15.98006579544212
$ws.Range("N20").Value = 23.11273885497573

# Row 21
$ws.Range("B21").Value = 16.5097097991255
$ws.Range("C21").Value = 6.050242212418357
$ws.Range("D21").Value = 7.975711657399073
$ws.Range("E21").Value = 10.1112312422008
$ws.Range("F21").Value = 39.57657362819626
$ws.Range("I21").Value = 32.18218140740146
$ws.Range("K21").Value = 14.31355007411588
$ws.Range("L21").Value = 10.38960536631309
$ws.Range("M21").Value = 16.02125457140689
$ws.Range("N21").Value = 23.05651782753882

# Row 22
$ws.Range("B22").Value = 16.67191692635516
$ws.Range("C22").Value = 6.237246793382288
$ws.Range("D22").Value = 7.971433868372952
$ws.Range("E22").Value = 10.1020034287545
$ws.Range("F22").Value = 39.67961902022457
$ws.Range("I22").Value = 32.18694372269988
$ws.Range("K22").Value = 14.43686192659659
$ws.Range("L22").Value = 10.38994311624116
$ws.Range("M22").Value = 16.0505498100231
$ws.Range("N22").Value = 23.02142755720844

# Row 23
$ws.Range("B23").Value = 16.58523796788649
$ws.Range("C23").Value = 6.138201672313364
$ws.Range("D23").Value = 7.973713321793771
$ws.Range("E23").Value = 10.1068643094124
$ws.Range("F23").Value = 39.62399203833271
$ws.Range("I23").Value = 32.18397384482434
$ws.Range("K23").Value = 14.37091885319929
$ws.Range("L23").Value = 10.38964660076755
$ws.Range("M23").Value = 16.03469873142124
$ws.Range("N23").Value = 23.04000470438888

# Row 24
$ws.Range("B24").Value = 16.25957459640244
$ws.Range("C24").Value = 5.74635230289848
$ws.Range("D24").Value = 7.982428115122421
$ws.Range("E24").Value = 10.12668602583533
$ws.Range("F24").Value = 39.42724157559473
$ws.Range("I24").Value = 32.18213531537744
$ws.Range("K24").Value = 14.12421770542448
$ws.Range("L24").Value = 10.39108160245739
$ws.Range("M24").Value = 15.97943307563026
$ws.Range("N24").Value = 23.11368569405304

# Row 25
$ws.Range("B25").Value = 15.91575153478305
$ws.Range("C25").Value = 5.288155719959961
$ws.Range("D25").Value = 7.992008220991409
$ws.Range("E25").Value = 10.1510960593542
$ws.Range("F25").Value = 39.24531947206527
$ws.Range("I25").Value = 32.2002661831935
$ws.Range("K25").Value = 13.86594494381699
$ws.Range("L25").Value = 10.39811795454702
$ws.Range("M25").Value = 15.93018514185314
$ws.Range("N25").Value = 23.20027190510187
